# Recreation of merges/pivot/chart/conditional-formatting rebuild that, for this
# particular sheet, manifests as:
#   1. Two header renames (player_id_x -> player_id, birth_year_x -> birth_year)
#   2. The 15 data rows (rows 2-16) being re-emitted in a different order
#      (matching the original row whose seas_id/calendar_year pairs are listed
#      below), with the player_id column (C) renumbered from 2193 to 3643.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames -----------------------------------------------------
$ws.Range("C1").Value = "player_id"
$ws.Range("E1").Value = "birth_year"

# --- 2. Re-order the 15 data rows ------------------------------------------
# Columns that hold text (everything else is numeric): B, D, F, I, J, AJ, AK, AL
$textCols = @("B", "D", "F", "I", "J", "AJ", "AK", "AL")
foreach ($col in $textCols) {
    $ws.Range($col + "2:" + $col + "16").NumberFormat = "@"
}

$dataRange = $ws.Range("A2:AM16")
$oldValues = $dataRange.Value2

# New row i (1-based, i = 1..15, representing sheet row i+1) takes its
# content from this original sheet row number.
$perm = @(8, 5, 10, 6, 3, 15, 14, 2, 9, 13, 7, 12, 16, 11, 4)

$newValues = New-Object 'object[,]' 15, 39
for ($i = 0; $i -lt 15; $i++) {
    $oldRowNum = $perm[$i]
    $oldIndex0 = $oldRowNum - 2
    for ($j = 0; $j -lt 39; $j++) {
        $newValues[$i, $j] = $oldValues[$oldIndex0 + 1, $j + 1]
    }
}
$dataRange.Value2 = $newValues

# --- 3. Renumber player_id column (C) from 2193 to 3643 -------------------
$ws.Range("C2:C16").Value = 3643
